$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: append/overwrite the two "movie list" readouts ---
# A2 already holds a string ("It") with a pre-existing direct format (quote-prefix
# style). Plain Value assignment on an existing string cell in this runtime drops
# that direct formatting, so refresh it via a copy/paste-values round trip through
# a scratch cell, which preserves the destination's existing style.
$ws.Range("ZZ1").Value = "Black Panther"
$ws.Range("ZZ1").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ1").ClearContents()

$ws.Range("A3").Value = "Interstellar"
$ws.Range("A4").Value = "It"
$ws.Range("A5").Value = "Star Wars: The Last Jedi"
$ws.Range("A6").Value = "Ready Player One"

# second movie list, appended directly under the first
$ws.Range("A7").Value = "Movie List 2"
$ws.Range("A8").Value = "Incredibles 2"
$ws.Range("A9").Value = "Avengers: Infinity War"
$ws.Range("A10").Value = "The Lego Batman Movie"
$ws.Range("A11").Value = "The Boss Baby"
$ws.Range("A12").Value = "Inside Out"

# --- Selection moves from B9 to A1:A2 ---
$ws.Range("A1:A2").Select() | Out-Null
